$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.56
$ws.Range("G2").Value = 1.7
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.96
$ws.Range("O2").Value = 1.3
$ws.Range("Q2").Value = 1.89
$ws.Range("R2").Value = 1.29
$ws.Range("S2").Value = 1.89
$ws.Range("T2").Value = 1.65
$ws.Range("U2").Value = 1.65
$ws.Range("V2").Value = 1.13
$ws.Range("W2").Value = 2.3
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 30
$ws.Range("Z2").Value = 75
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 11
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 34
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 13
$ws.Range("AG2").Value = 13.5
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 25
$ws.Range("AL2").Value = 50
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("G3").Value = 3.45
$ws.Range("J3").Value = 3.55
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 1.36
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 3.3
$ws.Range("O3").Value = 1.21
$ws.Range("R3").Value = 1.3
$ws.Range("S3").Value = 2.74
$ws.Range("T3").Value = 1.01
$ws.Range("U3").Value = 1.01
$ws.Range("V3").Value = 1.64
$ws.Range("W3").Value = 1.4
$ws.Range("X3").Value = 22
$ws.Range("Y3").Value = 15.5
$ws.Range("Z3").Value = 22
$ws.Range("AA3").Value = 44
$ws.Range("AB3").Value = 18.5
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 16
$ws.Range("AE3").Value = 36
$ws.Range("AF3").Value = 32
$ws.Range("AG3").Value = 19
$ws.Range("AH3").Value = 24
$ws.Range("AI3").Value = 50
$ws.Range("AJ3").Value = 75
$ws.Range("AK3").Value = 50
$ws.Range("AL3").Value = 60
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 5
$ws.Range("F5").Value = 2.94
$ws.Range("G5").Value = 3.65
$ws.Range("H5").Value = 2.24
$ws.Range("I5").Value = 2.68
$ws.Range("Q5").Value = 1.74

# Row 6
$ws.Range("G6").Value = 2.54

# Row 7
$ws.Range("F7").Value = 1.44
$ws.Range("H7").Value = 10
$ws.Range("Q7").Value = 2.32
$ws.Range("AK7").Value = 1000

# Row 10
$ws.Range("F10").Value = 1.48
$ws.Range("H10").Value = 1.04
$ws.Range("I10").Value = 20
$ws.Range("K10").Value = 980

# Row 11
$ws.Range("G11").Value = 2.48
$ws.Range("Q11").Value = 2.44
